$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell A1 from "Index" to "i"
$ws.Range("A1").Value = "i"

# Find the last used row
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Decrement every value in column A rows 2..lastRow by 1
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# Change column A width from 6 to 4 (stored width units)
$ws.Columns.Item(1).ColumnWidth = 3.17

Write-Host "Done"
